{"js": "const items = [\n  \"package com.example.myapplication\",\n  \"import android.os.Bundle\",\n  \"import android.widget.*\",\n  \"import androidx.appcompat.app.AppCompatActivity\",\n  \"import okhttp3.*\",\n  \"import okhttp3.MediaType.Companion.toMediaType\",\n  \"import okhttp3.RequestBody.Companion.toRequestBody\",\n  \"import org.json.JSONObject\",\n  \"import java.io.IOException\",\n  null,\n  \"class SecondActivity : AppCompatActivity() {\",\n  null,\n  \"    private val client = OkHttpClient()\",\n  null,\n  \"    override fun onCreate(savedInstanceState: Bundle?) {\",\n  \"        super.onCreate(savedInstanceState)\",\n  \"        setContentView(R.layout.activity_second)\",\n  null,\n  \"        val checkLunch = findViewById<CheckBox>(R.id.checkLunch)\",\n  \"        val checkDinner = findViewById<CheckBox>(R.id.checkDinner)\",\n  \"        val spinner = findViewById<Spinner>(R.id.spinnerLocation)\",\n  \"        val buttonSubmit = findViewById<Button>(R.id.buttonSubmit)\",\n  null,\n  \"        val locationList = listOf(\\\"V\u0103n Ph\u00f2ng MS\\\", \\\"Theo t\u1ed5 A Qu\u00fd\\\")\",\n  \"        val adapter = ArrayAdapter(this, android.R.layout.simple_spinner_item, locationList)\",\n  \"        adapter.setDropDownViewResource(android.R.layout.simple_spinner_dropdown_item)\",\n  \"        spinner.adapter = adapter\",\n  null,\n  \"        val msnv = intent.getStringExtra(\\\"msnv\\\") ?: \\\"unknown\\\"\",\n  null,\n  \"        buttonSubmit.setOnClickListener {\",\n  \"            val baocom = when {\",\n  \"                checkLunch.isChecked && checkDinner.isChecked -> \\\"tr\u01b0a,t\u1ed1i\\\"\",\n  \"                checkLunch.isChecked -> \\\"tr\u01b0a\\\"\",\n  \"                checkDinner.isChecked -> \\\"t\u1ed1i\\\"\",\n  \"                else -> {\",\n  \"                    Toast.makeText(this, \\\"Ch\u01b0a ch\u1ecdn b\u00e1o c\u01a1m n\u00e0o\\\", Toast.LENGTH_SHORT).show()\",\n  \"                    return@setOnClickListener\",\n  \"                }\",\n  \"            }\",\n  null,\n  \"            val vitri = spinner.selectedItem.toString()\",\n  \"            sendBaoCom(msnv, baocom, vitri)\",\n  \"        }\",\n  \"    }\",\n  null,\n  \"    private fun sendBaoCom(msnv: String, baocom: String, vitri: String) {\",\n  \"        val json = JSONObject().apply {\",\n  \"            put(\\\"msnv\\\", msnv)\",\n  \"            put(\\\"baocom\\\", baocom)\",\n  \"            put(\\\"vitri\\\", vitri)\",\n  \"        }\",\n  null,\n  \"        val requestBody = json.toString()\",\n  \"            .toRequestBody(\\\"application/json\\\".toMediaType())\",\n  null,\n  \"        val request = Request.Builder()\",\n  \"            .url(\\\"http://192.168.1.100:5000/baocom\\\") // thay \u0111\u00fang IP m\u00e1y Flask\",\n  \"            .post(requestBody)\",\n  \"            .build()\",\n  null,\n  \"        client.newCall(request).enqueue(object : Callback {\",\n  \"            override fun onFailure(call: Call, e: IOException) {\",\n  \"                runOnUiThread {\",\n  \"                    Toast.makeText(applicationContext, \\\"L\u1ed7i g\u1eedi b\u00e1o c\u01a1m\\\", Toast.LENGTH_SHORT).show()\",\n  \"                }\",\n  \"            }\",\n  null,\n  \"            override fun onResponse(call: Call, response: Response) {\",\n  \"                runOnUiThread {\",\n  \"                    val message = if (response.isSuccessful) \\\"\u0110\u00e3 g\u1eedi b\u00e1o c\u01a1m!\\\" else \\\"G\u1eedi th\u1ea5t b\u1ea1i!\\\"\",\n  \"                    Toast.makeText(applicationContext, message, Toast.LENGTH_SHORT).show()\",\n  \"                }\",\n  \"            }\",\n  \"        })\",\n  \"    }\",\n  \"}\",\n  null,\n  \"PQ-KLBT-STR-DTL-PMC-01004-01-AS-0006\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet current = paragraphs.items[paragraphs.items.length - 1];\n\nfor (const text of items) {\n  current = current.insertParagraph(text === null ? \"\" : text, Word.InsertLocation.after);\n  await context.sync();\n}\n", "ps1": "$items = @(\n    'package com.example.myapplication'\n    'import android.os.Bundle'\n    'import android.widget.*'\n    'import androidx.appcompat.app.AppCompatActivity'\n    'import okhttp3.*'\n    'import okhttp3.MediaType.Companion.toMediaType'\n    'import okhttp3.RequestBody.Companion.toRequestBody'\n    'import org.json.JSONObject'\n    'import java.io.IOException'\n    ''\n    'class SecondActivity : AppCompatActivity() {'\n    ''\n    '    private val client = OkHttpClient()'\n    ''\n    '    override fun onCreate(savedInstanceState: Bundle?) {'\n    '        super.onCreate(savedInstanceState)'\n    '        setContentView(R.layout.activity_second)'\n    ''\n    '        val checkLunch = findViewById<CheckBox>(R.id.checkLunch)'\n    '        val checkDinner = findViewById<CheckBox>(R.id.checkDinner)'\n    '        val spinner = findViewById<Spinner>(R.id.spinnerLocation)'\n    '        val buttonSubmit = findViewById<Button>(R.id.buttonSubmit)'\n    ''\n    '        val locationList = listOf(\"V\u0103n Ph\u00f2ng MS\", \"Theo t\u1ed5 A Qu\u00fd\")'\n    '        val adapter = ArrayAdapter(this, android.R.layout.simple_spinner_item, locationList)'\n    '        adapter.setDropDownViewResource(android.R.layout.simple_spinner_dropdown_item)'\n    '        spinner.adapter = adapter'\n    ''\n    '        val msnv = intent.getStringExtra(\"msnv\") ?: \"unknown\"'\n    ''\n    '        buttonSubmit.setOnClickListener {'\n    '            val baocom = when {'\n    '                checkLunch.isChecked && checkDinner.isChecked -> \"tr\u01b0a,t\u1ed1i\"'\n    '                checkLunch.isChecked -> \"tr\u01b0a\"'\n    '                checkDinner.isChecked -> \"t\u1ed1i\"'\n    '                else -> {'\n    '                    Toast.makeText(this, \"Ch\u01b0a ch\u1ecdn b\u00e1o c\u01a1m n\u00e0o\", Toast.LENGTH_SHORT).show()'\n    '                    return@setOnClickListener'\n    '                }'\n    '            }'\n    ''\n    '            val vitri = spinner.selectedItem.toString()'\n    '            sendBaoCom(msnv, baocom, vitri)'\n    '        }'\n    '    }'\n    ''\n    '    private fun sendBaoCom(msnv: String, baocom: String, vitri: String) {'\n    '        val json = JSONObject().apply {'\n    '            put(\"msnv\", msnv)'\n    '            put(\"baocom\", baocom)'\n    '            put(\"vitri\", vitri)'\n    '        }'\n    ''\n    '        val requestBody = json.toString()'\n    '            .toRequestBody(\"application/json\".toMediaType())'\n    ''\n    '        val request = Request.Builder()'\n    '            .url(\"http://192.168.1.100:5000/baocom\") // thay \u0111\u00fang IP m\u00e1y Flask'\n    '            .post(requestBody)'\n    '            .build()'\n    ''\n    '        client.newCall(request).enqueue(object : Callback {'\n    '            override fun onFailure(call: Call, e: IOException) {'\n    '                runOnUiThread {'\n    '                    Toast.makeText(applicationContext, \"L\u1ed7i g\u1eedi b\u00e1o c\u01a1m\", Toast.LENGTH_SHORT).show()'\n    '                }'\n    '            }'\n    ''\n    '            override fun onResponse(call: Call, response: Response) {'\n    '                runOnUiThread {'\n    '                    val message = if (response.isSuccessful) \"\u0110\u00e3 g\u1eedi b\u00e1o c\u01a1m!\" else \"G\u1eedi th\u1ea5t b\u1ea1i!\"'\n    '                    Toast.makeText(applicationContext, message, Toast.LENGTH_SHORT).show()'\n    '                }'\n    '            }'\n    '        })'\n    '    }'\n    '}'\n    ''\n    'PQ-KLBT-STR-DTL-PMC-01004-01-AS-0006'\n)\n\n$d = $word.ActiveDocument\n\nforeach ($t in $items) {\n    $last = $d.Paragraphs.Last\n    $last.Range.InsertParagraphAfter()\n    if ($t -ne '') {\n        $d.Paragraphs.Last.Range.Text = $t\n    }\n}\n"}
